$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting (styles, fill, borders, number format) from row 40 into row 41
$ws.Range("A40:C40").Copy()
$ws.Range("A41:C41").PasteSpecial(-4122)

# Set the new row's values
$ws.Range("A41").Value = 43349
$ws.Range("B41").Value = "python/pip/virtualenv"
$ws.Range("C41").Value = "installed virtualenv and virtualenvwrapper ,created a virtualenv"

$ws.Range("C41").Select()
